{"js": "const searchResults = context.document.body.search(\" souhv\u011bzd\u00ed Souhv\u011bzd\u00ed\", { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < searchResults.items.length; i++) {\n  searchResults.items[i].insertText(\"Souhv\u011bzd\u00ed\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$find.Text = \" souhv\u011bzd\u00ed Souhv\u011bzd\u00ed\"\n$find.Replacement.Text = \"Souhv\u011bzd\u00ed\"\n\n# wdFindContinue=1, wdReplaceAll=2\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
